# Update countries & provincias Spain
# - refresh "Datos actualizados" timestamp
# - update per-country case/recovered/death statistics
# - re-sort a handful of rows whose "Casos totales" changed enough to
#   move them relative to their neighbours (sheet is kept sorted
#   descending by column B), which surfaces as country-name swaps on
#   the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = 'Datos actualizados a 23 de Julio de 2020 a las 13:22'

# Iran
$ws.Range("B14").Value = 284034
$ws.Range("C14").Value = 2621
$ws.Range("D14").Value = 247230
$ws.Range("E14").Value = 21730
$ws.Range("G14").Value = 221
$ws.Range("H14").Value = 15074

# Emiratos Arabes Unidos
$ws.Range("B40").Value = 57988
$ws.Range("C40").Value = 254
$ws.Range("D40").Value = 50848
$ws.Range("E40").Value = 6798

# Suiza
$ws.Range("B55").Value = 34000
$ws.Range("C55").Value = 117
$ws.Range("E55").Value = 1528

# Nepal
$ws.Range("B66").Value = 18241
$ws.Range("C66").Value = 147
$ws.Range("D66").Value = 12840
$ws.Range("E66").Value = 5358
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 43

# Rows 80-84: Estado de Palestina moves up ahead of Bulgaria / Republica
# de Macedonia, and Bosnia y Herzegovina moves up ahead of Senegal
$ws.Range("A80").Value = 'Estado de Palestina'
$ws.Range("B80").Value = 9744
$ws.Range("C80").Value = 346
$ws.Range("D80").Value = 2720
$ws.Range("E80").Value = 6958
$ws.Range("H80").Value = 66

$ws.Range("A81").Value = 'Bulgaria'
$ws.Range("B81").Value = 9584
$ws.Range("D81").Value = 4643
$ws.Range("E81").Value = 4620
$ws.Range("H81").Value = 321

$ws.Range("A82").Value = 'Republica de Macedonia'
$ws.Range("B82").Value = 9547
$ws.Range("D82").Value = 5071
$ws.Range("E82").Value = 4034
$ws.Range("H82").Value = 442

$ws.Range("A83").Value = 'Bosnia y Herzegovina'
$ws.Range("B83").Value = 9462
$ws.Range("C83").Value = 347
$ws.Range("D83").Value = 4367
$ws.Range("E83").Value = 4821
$ws.Range("G83").Value = 6
$ws.Range("H83").Value = 274

$ws.Range("A84").Value = 'Senegal'
$ws.Range("B84").Value = 9266
$ws.Range("C84").Value = 145
$ws.Range("D84").Value = 6170
$ws.Range("E84").Value = 2918
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 178

# Consejo Danes para los Refugiados
$ws.Range("B87").Value = 8720
$ws.Range("C87").Value = 94
$ws.Range("D87").Value = 5105
$ws.Range("E87").Value = 3414
$ws.Range("G87").Value = 4
$ws.Range("H87").Value = 201

# Madagascar
$ws.Range("B88").Value = 8381
$ws.Range("C88").Value = 219
$ws.Range("D88").Value = 5160
$ws.Range("E88").Value = 3151
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 70

# Row 156 (unaffected country identity)
$ws.Range("B156").Value = 680
$ws.Range("C156").Value = 1
$ws.Range("E156").Value = 6

# Rows 165-169: Burundi moves up ahead of Birmania / Mauricio / Comoras
# / Isla de Man
$ws.Range("A165").Value = 'Burundi'
$ws.Range("B165").Value = 345
$ws.Range("C165").Value = 17
$ws.Range("D165").Value = 270
$ws.Range("E165").Value = 74
$ws.Range("H165").Value = 1

$ws.Range("A166").Value = 'Birmania'
$ws.Range("D166").Value = 280
$ws.Range("E166").Value = 57
$ws.Range("H166").Value = 6

$ws.Range("A167").Value = 'Mauricio'
$ws.Range("B167").Value = 343
$ws.Range("D167").Value = 332
$ws.Range("E167").Value = 1
$ws.Range("H167").Value = 10

$ws.Range("A168").Value = 'Comoras'
$ws.Range("B168").Value = 337
$ws.Range("D168").Value = 319
$ws.Range("E168").Value = 11
$ws.Range("H168").Value = 7

$ws.Range("A169").Value = 'Isla de Man'
$ws.Range("B169").Value = 336
$ws.Range("D169").Value = 312
$ws.Range("E169").Value = 0
$ws.Range("H169").Value = 24

# Rows 210-211: Groenlandia moves up ahead of Islas Malvinas
$ws.Range("A210").Value = 'Groenlandia'
$ws.Range("A211").Value = 'Islas Malvinas'
